$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Rows.Item(4).Delete()
